$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.3269336956678857
$ws.Range("J2").Value = 0.3269336956678857
$ws.Range("M2").Value = 0.6189250000000001
$ws.Range("N2").Value = 1.856775
$ws.Range("O2").Value = 0.09614699503454774
$ws.Range("P2").Value = 0.09614699503454775
$ws.Range("Q2").Value = 0.06999299040000001
$ws.Range("R2").Value = 0.6299369136
$ws.Range("S2").Value = 0.03143369241400655
$ws.Range("T2").Value = 0.03143369241400656
$ws.Range("I3").Value = 0.3269336956678857
$ws.Range("J3").Value = 0.3269336956678857
$ws.Range("O3").Value = 0.7879294335349575
$ws.Range("P3").Value = 0.7879294335349576
$ws.Range("S3").Value = 0.2576006816310873
$ws.Range("T3").Value = 0.2576006816310874
$ws.Range("I4").Value = 0.3269336956678857
$ws.Range("J4").Value = 0.3269336956678857
$ws.Range("M4").Value = 0.01220666666666667
$ws.Range("N4").Value = 0.03662
$ws.Range("O4").Value = 0.001896246426284896
$ws.Range("P4").Value = 0.001896246426284896
$ws.Range("Q4").Value = 0.00138042752
$ws.Range("R4").Value = 0.01242384768
$ws.Range("S4").Value = 0.0006199468520423421
$ws.Range("T4").Value = 0.0006199468520423422
$ws.Range("I5").Value = 0.3269336956678857
$ws.Range("J5").Value = 0.3269336956678857
$ws.Range("M5").Value = 0.7340256666666667
$ws.Range("N5").Value = 2.202077
$ws.Range("O5").Value = 0.1140273250042099
$ws.Range("P5").Value = 0.1140273250042099
$ws.Range("Q5").Value = 0.08300949459200001
$ws.Range("R5").Value = 0.7470854513280001
$ws.Range("S5").Value = 0.03727937477074945
$ws.Range("T5").Value = 0.03727937477074945
$ws.Range("A6").Value = 'Inflammatory-Mac'
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.084206
$ws.Range("H6").Value = 0.252618
$ws.Range("I6").Value = 0.2434367817753429
$ws.Range("J6").Value = 0.243436781775343
$ws.Range("M6").Value = 0.6189250000000001
$ws.Range("N6").Value = 1.856775
$ws.Range("O6").Value = 0.09614699503454774
$ws.Range("P6").Value = 0.09614699503454775
$ws.Range("Q6").Value = 0.05211719855000001
$ws.Range("R6").Value = 0.46905478695
$ws.Range("S6").Value = 0.02340571504858018
$ws.Range("T6").Value = 0.02340571504858019
$ws.Range("A7").Value = 'Inflammatory-Mac'
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.084206
$ws.Range("H7").Value = 0.252618
$ws.Range("I7").Value = 0.2434367817753429
$ws.Range("J7").Value = 0.243436781775343
$ws.Range("O7").Value = 0.7879294335349575
$ws.Range("P7").Value = 0.7879294335349576
$ws.Range("Q7").Value = 0.4271030489946666
$ws.Range("R7").Value = 3.843927440952
$ws.Range("S7").Value = 0.191811005565819
$ws.Range("T7").Value = 0.1918110055658191
$ws.Range("A8").Value = 'Inflammatory-Mac'
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.3333333333333333
$ws.Range("G8").Value = 0.084206
$ws.Range("H8").Value = 0.252618
$ws.Range("I8").Value = 0.2434367817753429
$ws.Range("J8").Value = 0.243436781775343
$ws.Range("M8").Value = 0.01220666666666667
$ws.Range("N8").Value = 0.03662
$ws.Range("O8").Value = 0.001896246426284896
$ws.Range("P8").Value = 0.001896246426284896
$ws.Range("Q8").Value = 0.001027874573333333
$ws.Range("R8").Value = 0.00925087116
$ws.Range("S8").Value = 0.0004616161274677902
$ws.Range("T8").Value = 0.0004616161274677903
$ws.Range("A9").Value = 'Inflammatory-Mac'
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.3333333333333333
$ws.Range("G9").Value = 0.084206
$ws.Range("H9").Value = 0.252618
$ws.Range("I9").Value = 0.2434367817753429
$ws.Range("J9").Value = 0.243436781775343
$ws.Range("M9").Value = 0.7340256666666667
$ws.Range("N9").Value = 2.202077
$ws.Range("O9").Value = 0.1140273250042099
$ws.Range("P9").Value = 0.1140273250042099
$ws.Range("Q9").Value = 0.06180936528733334
$ws.Range("R9").Value = 0.556284287586
$ws.Range("S9").Value = 0.02775844503347594
$ws.Range("T9").Value = 0.02775844503347595
$ws.Range("A10").Value = 'MuSCs'
$ws.Range("G10").Value = 0.08859233333333333
$ws.Range("H10").Value = 0.265777
$ws.Range("I10").Value = 0.2561175274521424
$ws.Range("J10").Value = 0.2561175274521424
$ws.Range("M10").Value = 0.6189250000000001
$ws.Range("N10").Value = 1.856775
$ws.Range("O10").Value = 0.09614699503454774
$ws.Range("P10").Value = 0.09614699503454775
$ws.Range("Q10").Value = 0.05483200990833333
$ws.Range("R10").Value = 0.493488089175
$ws.Range("S10").Value = 0.02462493064020178
$ws.Range("T10").Value = 0.02462493064020179
$ws.Range("A11").Value = 'MuSCs'
$ws.Range("G11").Value = 0.08859233333333333
$ws.Range("H11").Value = 0.265777
$ws.Range("I11").Value = 0.2561175274521424
$ws.Range("J11").Value = 0.2561175274521424
$ws.Range("O11").Value = 0.7879294335349575
$ws.Range("P11").Value = 0.7879294335349576
$ws.Range("Q11").Value = 0.4493510638697777
$ws.Range("R11").Value = 4.044159574827999
$ws.Range("S11").Value = 0.2018025383237405
$ws.Range("T11").Value = 0.2018025383237405
$ws.Range("A12").Value = 'MuSCs'
$ws.Range("G12").Value = 0.08859233333333333
$ws.Range("H12").Value = 0.265777
$ws.Range("I12").Value = 0.2561175274521424
$ws.Range("J12").Value = 0.2561175274521424
$ws.Range("M12").Value = 0.01220666666666667
$ws.Range("N12").Value = 0.03662
$ws.Range("O12").Value = 0.001896246426284896
$ws.Range("P12").Value = 0.001896246426284896
$ws.Range("Q12").Value = 0.001081417082222222
$ws.Range("R12").Value = 0.009732753739999999
$ws.Range("S12").Value = 0.0004856619461400489
$ws.Range("T12").Value = 0.0004856619461400489
$ws.Range("A13").Value = 'MuSCs'
$ws.Range("G13").Value = 0.08859233333333333
$ws.Range("H13").Value = 0.265777
$ws.Range("I13").Value = 0.2561175274521424
$ws.Range("J13").Value = 0.2561175274521424
$ws.Range("M13").Value = 0.7340256666666667
$ws.Range("N13").Value = 2.202077
$ws.Range("O13").Value = 0.1140273250042099
$ws.Range("P13").Value = 0.1140273250042099
$ws.Range("Q13").Value = 0.06502904653655556
$ws.Range("R13").Value = 0.585261418829
$ws.Range("S13").Value = 0.02920439654206009
$ws.Range("T13").Value = 0.02920439654206009
$ws.Range("G14").Value = 0.06001866666666666
$ws.Range("H14").Value = 0.180056
$ws.Range("I14").Value = 0.1735119951046289
$ws.Range("J14").Value = 0.1735119951046289
$ws.Range("M14").Value = 0.6189250000000001
$ws.Range("N14").Value = 1.856775
$ws.Range("O14").Value = 0.09614699503454774
$ws.Range("P14").Value = 0.09614699503454775
$ws.Range("Q14").Value = 0.03714705326666667
$ws.Range("R14").Value = 0.3343234794
$ws.Range("S14").Value = 0.01668265693175923
$ws.Range("T14").Value = 0.01668265693175923
$ws.Range("G15").Value = 0.06001866666666666
$ws.Range("H15").Value = 0.180056
$ws.Range("I15").Value = 0.1735119951046289
$ws.Range("J15").Value = 0.1735119951046289
$ws.Range("O15").Value = 0.7879294335349575
$ws.Range("P15").Value = 0.7879294335349576
$ws.Range("Q15").Value = 0.3044219595982222
$ws.Range("R15").Value = 2.739797636384
$ws.Range("S15").Value = 0.1367152080143106
$ws.Range("T15").Value = 0.1367152080143106
$ws.Range("G16").Value = 0.06001866666666666
$ws.Range("H16").Value = 0.180056
$ws.Range("I16").Value = 0.1735119951046289
$ws.Range("J16").Value = 0.1735119951046289
$ws.Range("M16").Value = 0.01220666666666667
$ws.Range("N16").Value = 0.03662
$ws.Range("O16").Value = 0.001896246426284896
$ws.Range("P16").Value = 0.001896246426284896
$ws.Range("Q16").Value = 0.0007326278577777778
$ws.Range("R16").Value = 0.00659365072
$ws.Range("S16").Value = 0.000329021500634715
$ws.Range("T16").Value = 0.000329021500634715
$ws.Range("G17").Value = 0.06001866666666666
$ws.Range("H17").Value = 0.180056
$ws.Range("I17").Value = 0.1735119951046289
$ws.Range("J17").Value = 0.1735119951046289
$ws.Range("M17").Value = 0.7340256666666667
$ws.Range("N17").Value = 2.202077
$ws.Range("O17").Value = 0.1140273250042099
$ws.Range("P17").Value = 0.1140273250042099
$ws.Range("Q17").Value = 0.04405524181244445
$ws.Range("R17").Value = 0.396497176312
$ws.Range("S17").Value = 0.01978510865792439
$ws.Range("T17").Value = 0.01978510865792439